$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.919190342406841
$ws.Range("D2").Value = 4.492198819804734
$ws.Range("E2").Value = 11.10450525908647
$ws.Range("F2").Value = 24.03736573980291
$ws.Range("G2").Value = 3.620191562047249
$ws.Range("K2").Value = 12.25411080649266
$ws.Range("M2").Value = 14.95469571531774
$ws.Range("N2").Value = 17.10590391416503
$ws.Range("O2").Value = 21.33100859478309
$ws.Range("C3").Value = 4.749640844105196
$ws.Range("D3").Value = 4.467701378385796
$ws.Range("E3").Value = 10.99872644801263
$ws.Range("F3").Value = 23.97770994788492
$ws.Range("G3").Value = 3.62250863713381
$ws.Range("K3").Value = 11.78661070194214
$ws.Range("M3").Value = 14.67708806620432
$ws.Range("N3").Value = 17.16135385502553
$ws.Range("O3").Value = 21.35098937871329
$ws.Range("C4").Value = 4.643904877941972
$ws.Range("D4").Value = 4.452439006721589
$ws.Range("E4").Value = 10.93767655890509
$ws.Range("F4").Value = 23.94915712596292
$ws.Range("G4").Value = 3.62400640022239
$ws.Range("K4").Value = 11.49145037029118
$ws.Range("M4").Value = 14.50747125100855
$ws.Range("N4").Value = 17.19724032626039
$ws.Range("O4").Value = 21.3696132621246
$ws.Range("C5").Value = 4.600483853942342
$ws.Range("D5").Value = 4.446166120144294
$ws.Range("E5").Value = 10.91380462337898
$ws.Range("F5").Value = 23.93955875734424
$ws.Range("G5").Value = 3.624635690253007
$ws.Range("K5").Value = 11.36931213313382
$ws.Range("M5").Value = 14.43865730251306
$ws.Range("N5").Value = 17.21232812424528
$ws.Range("O5").Value = 21.37879584239701
$ws.Range("C6").Value = 4.593256115459659
$ws.Range("D6").Value = 4.445121354032493
$ws.Range("E6").Value = 10.90990220079134
$ws.Range("F6").Value = 23.93808812715315
$ws.Range("G6").Value = 3.624741329115984
$ws.Range("K6").Value = 11.34892466140741
$ws.Range("M6").Value = 14.42725212546023
$ws.Range("N6").Value = 17.21486148454236
$ws.Range("O6").Value = 21.38041666481251
$ws.Range("C7").Value = 4.643320527103963
$ws.Range("D7").Value = 4.452354621501254
$ws.Range("E7").Value = 10.93735050630575
$ws.Range("F7").Value = 23.94901942435161
$ws.Range("G7").Value = 3.62401481028059
$ws.Range("K7").Value = 11.48981043921051
$ws.Range("M7").Value = 14.50654183417641
$ws.Range("N7").Value = 17.19744192635295
$ws.Range("O7").Value = 21.36973065811992
$ws.Range("C8").Value = 4.861113011337875
$ws.Range("D8").Value = 4.483799246555555
$ws.Range("E8").Value = 11.06723991962503
$ws.Range("F8").Value = 24.0151248357711
$ws.Range("G8").Value = 3.620974946549825
$ws.Range("K8").Value = 12.09470398726273
$ws.Range("M8").Value = 14.85886060418095
$ws.Range("N8").Value = 17.12464181934358
$ws.Range("O8").Value = 21.33657614632773
$ws.Range("C9").Value = 5.272163969915333
$ws.Range("D9").Value = 4.543615083202591
$ws.Range("E9").Value = 11.35162524044515
$ws.Range("F9").Value = 24.2084361706275
$ws.Range("G9").Value = 3.61560657214342
$ws.Range("K9").Value = 13.20949019163055
$ws.Range("M9").Value = 15.55222675023141
$ws.Range("N9").Value = 16.99643106993837
$ws.Range("O9").Value = 21.32217094293364
$ws.Range("C10").Value = 5.560695030314664
$ws.Range("D10").Value = 4.586311053666872
$ws.Range("E10").Value = 11.57681499954575
$ws.Range("F10").Value = 24.38856938935929
$ws.Range("G10").Value = 3.612019810874936
$ws.Range("K10").Value = 13.97685141343373
$ws.Range("M10").Value = 16.057734956447
$ws.Range("N10").Value = 16.91103612584481
$ws.Range("O10").Value = 21.34262951477648
$ws.Range("C11").Value = 5.688376376451021
$ws.Range("D11").Value = 4.605436710552919
$ws.Range("E11").Value = 11.68236337076397
$ws.Range("F11").Value = 24.47858530007077
$ws.Range("G11").Value = 3.610464851669803
$ws.Range("K11").Value = 14.31333835451626
$ws.Range("M11").Value = 16.2857696544581
$ws.Range("N11").Value = 16.87408468600064
$ws.Range("O11").Value = 21.35869622154159
$ws.Range("C12").Value = 5.736165685089585
$ws.Range("D12").Value = 4.612634163387087
$ws.Range("E12").Value = 11.72274149405374
$ws.Range("F12").Value = 24.51381157892769
$ws.Range("G12").Value = 3.609886990366725
$ws.Range("K12").Value = 14.43885120568859
$ws.Range("M12").Value = 16.37175859021953
$ws.Range("N12").Value = 16.86036363540807
$ws.Range("O12").Value = 21.36575208017066
$ws.Range("C13").Value = 5.725899037290141
$ws.Range("D13").Value = 4.611086102001356
$ws.Range("E13").Value = 11.71402774963156
$ws.Range("F13").Value = 24.50617467161437
$ws.Range("G13").Value = 3.610010956251403
$ws.Range("K13").Value = 14.41190600819308
$ws.Range("M13").Value = 16.35325679852124
$ws.Range("N13").Value = 16.86330664054263
$ws.Range("O13").Value = 21.36418927317379
$ws.Range("C14").Value = 5.692319521869635
$ws.Range("D14").Value = 4.606029764759489
$ws.Range("E14").Value = 11.6856773200186
$ws.Range("F14").Value = 24.48146068783195
$ws.Range("G14").Value = 3.610417091188098
$ws.Range("K14").Value = 14.32370310977156
$ws.Range("M14").Value = 16.29285176639052
$ws.Range("N14").Value = 16.87295040732916
$ws.Range("O14").Value = 21.35925724188664
$ws.Range("C15").Value = 5.67167675761683
$ws.Range("D15").Value = 4.60292668108201
$ws.Range("E15").Value = 11.66836402695403
$ws.Range("F15").Value = 24.4664703366279
$ws.Range("G15").Value = 3.610667287207165
$ws.Range("K15").Value = 14.26942522553713
$ws.Range("M15").Value = 16.25580220122743
$ws.Range("N15").Value = 16.87889284566552
$ws.Range("O15").Value = 21.35636275151903
$ws.Range("C16").Value = 5.552274707752786
$ws.Range("D16").Value = 4.585054973028118
$ws.Range("E16").Value = 11.5699765056471
$ws.Range("F16").Value = 24.38284728330247
$ws.Range("G16").Value = 3.612122969076978
$ws.Range("K16").Value = 13.95459926146707
$ws.Range("M16").Value = 16.04278617296535
$ws.Range("N16").Value = 16.91348905024089
$ws.Range("O16").Value = 21.34171554275282
$ws.Range("C17").Value = 5.478075046220011
$ws.Range("D17").Value = 4.574013547264542
$ws.Range("E17").Value = 11.51038848834273
$ws.Range("F17").Value = 24.33360017595515
$ws.Range("G17").Value = 3.613035580055688
$ws.Range("K17").Value = 13.75816587072919
$ws.Range("M17").Value = 15.91155232748874
$ws.Range("N17").Value = 16.9351974720736
$ws.Range("O17").Value = 21.33446146426643
$ws.Range("C18").Value = 5.435063513641696
$ws.Range("D18").Value = 4.567634997818746
$ws.Range("E18").Value = 11.47640993457667
$ws.Range("F18").Value = 24.30603580234739
$ws.Range("G18").Value = 3.613567710346195
$ws.Range("K18").Value = 13.64400249746749
$ws.Range("M18").Value = 15.8358939316565
$ws.Range("N18").Value = 16.94786201227169
$ws.Range("O18").Value = 21.33092548288544
$ws.Range("C19").Value = 5.420444734511689
$ws.Range("D19").Value = 4.565470624530802
$ws.Range("E19").Value = 11.46495716632377
$ws.Range("F19").Value = 24.29683434539195
$ws.Range("G19").Value = 3.613749122436683
$ws.Range("K19").Value = 13.6051492396247
$ws.Range("M19").Value = 15.81024974314901
$ws.Range("N19").Value = 16.95218068026429
$ws.Range("O19").Value = 21.32983755018243
$ws.Range("C20").Value = 5.4860086376521
$ws.Range("D20").Value = 4.575191821301861
$ws.Range("E20").Value = 11.51670148473972
$ws.Range("F20").Value = 24.33876397594513
$ws.Range("G20").Value = 3.612937684175579
$ws.Range("K20").Value = 13.77919942017384
$ws.Range("M20").Value = 15.92554126493778
$ws.Range("N20").Value = 16.93286811297391
$ws.Range("O20").Value = 21.33516780999707
$ws.Range("C21").Value = 5.702198208564433
$ws.Range("D21").Value = 4.607516173717153
$ws.Range("E21").Value = 11.69399372583026
$ws.Range("F21").Value = 24.48868904988172
$ws.Range("G21").Value = 3.610297502303847
$ws.Range("K21").Value = 14.34966291521644
$ws.Range("M21").Value = 16.31060470174905
$ws.Range("N21").Value = 16.87011043328019
$ws.Range("O21").Value = 21.36067953426487
$ws.Range("C22").Value = 5.840203519346346
$ws.Range("D22").Value = 4.628378546235757
$ws.Range("E22").Value = 11.81223192488783
$ws.Range("F22").Value = 24.59330228570424
$ws.Range("G22").Value = 3.608635894109601
$ws.Range("K22").Value = 14.71133969253713
$ws.Range("M22").Value = 16.56011409741731
$ws.Range("N22").Value = 16.83067770988318
$ws.Range("O22").Value = 21.3830160624081
$ws.Range("C23").Value = 5.766862400872416
$ws.Range("D23").Value = 4.617268771785588
$ws.Range("E23").Value = 11.74892198158065
$ws.Range("F23").Value = 24.53686945340028
$ws.Range("G23").Value = 3.609516897573751
$ws.Range("K23").Value = 14.51935532059529
$ws.Range("M23").Value = 16.42717046530245
$ws.Range("N23").Value = 16.85157912539525
$ws.Range("O23").Value = 21.37057685481727
$ws.Range("C24").Value = 5.482422957377173
$ws.Range("D24").Value = 4.574659218997339
$ws.Range("E24").Value = 11.51384650594158
$ws.Range("F24").Value = 24.33642708894383
$ws.Range("G24").Value = 3.612981919666668
$ws.Range("K24").Value = 13.76969397986551
$ws.Range("M24").Value = 15.91921751475781
$ws.Range("N24").Value = 16.93392064272653
$ws.Range("O24").Value = 21.33484649463463
$ws.Range("C25").Value = 5.163086559012145
$ws.Range("D25").Value = 4.527643954848259
$ws.Range("E25").Value = 11.27169744463145
$ws.Range("F25").Value = 24.14938631533024
$ws.Range("G25").Value = 3.616995812161371
$ws.Range("K25").Value = 12.9164631402916
$ws.Range("M25").Value = 15.36495842309112
$ws.Range("N25").Value = 17.21232812424528
$ws.Range("O25").Value = 21.32062622070966
